$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 64, pushing the existing rows 64-80 down to 65-81.
$ws.Rows(64).Insert()

# Populate the newly inserted row 64 with the new record.
$ws.Cells.Item(64, 1).Value  = 5
$ws.Cells.Item(64, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(64, 3).Value  = "Maule"
$ws.Cells.Item(64, 4).Value  = 44876
$ws.Cells.Item(64, 5).Value  = 7
$ws.Cells.Item(64, 6).Value  = 100112040
$ws.Cells.Item(64, 7).Value  = "Cilantro"
$ws.Cells.Item(64, 8).Value  = "Sin especificar"
$ws.Cells.Item(64, 9).Value  = "Primera"
$ws.Cells.Item(64, 10).Value = 150
$ws.Cells.Item(64, 11).Value = 7000
$ws.Cells.Item(64, 12).Value = 7000
$ws.Cells.Item(64, 13).Value = 7000
$ws.Cells.Item(64, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(64, 15).Value = "Región del Maule"
$ws.Cells.Item(64, 16).Value = 194
$ws.Cells.Item(64, 17).Value = 36
$ws.Cells.Item(64, 18).Value = "Hortaliza"
